$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("H2").Value = 7242
$ws.Range("I2").Value = 7279
$ws.Range("J2").Value = 4829
$ws.Range("H3").Value = 8347
$ws.Range("I3").Value = 7488
$ws.Range("J3").Value = 5118
$ws.Range("I4").Value = 1773
$ws.Range("J4").Value = 1141
$ws.Range("J5").Value = 406
$ws.Range("J6").Value = 6330
$ws.Range("I7").Value = 26223
$ws.Range("J7").Value = 17824

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 202
$ws.Range("J3").Value = 283
$ws.Range("J7").Value = 778

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J3").Value = 90
$ws.Range("J5").Value = 8
$ws.Range("J6").Value = 65
$ws.Range("J7").Value = 269

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J2").Value = 75
$ws.Range("J7").Value = 266

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J3").Value = 266
$ws.Range("J4").Value = 51
$ws.Range("J5").Value = 19
$ws.Range("J6").Value = 183
$ws.Range("J7").Value = 687

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J2").Value = 60
$ws.Range("J7").Value = 162

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("J2").Value = 49
$ws.Range("J7").Value = 149

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J2").Value = 135
$ws.Range("J3").Value = 7
$ws.Range("J6").Value = 133
$ws.Range("J7").Value = 515
$ws.Range("J8").Value = 1141
$ws.Range("J10").Value = 118
$ws.Range("J11").Value = 272
$ws.Range("J13").Value = 21
$ws.Range("J14").Value = 76
$ws.Range("J17").Value = 26
$ws.Range("J19").Value = 520
$ws.Range("J20").Value = 373
$ws.Range("J23").Value = 169
$ws.Range("J27").Value = 96
$ws.Range("J29").Value = 1012
$ws.Range("J31").Value = 162
$ws.Range("J32").Value = 28
$ws.Range("J33").Value = 812
$ws.Range("J35").Value = 28
$ws.Range("J36").Value = 247
$ws.Range("J37").Value = 557
$ws.Range("J41").Value = 115
$ws.Range("J42").Value = 722
$ws.Range("J44").Value = 133
$ws.Range("J48").Value = 199
$ws.Range("J50").Value = 105
$ws.Range("J51").Value = 223
$ws.Range("J52").Value = 451
$ws.Range("I63").Value = 234
$ws.Range("J63").Value = 65
$ws.Range("J64").Value = 119
$ws.Range("J66").Value = 57
$ws.Range("J67").Value = 687
$ws.Range("J73").Value = 164
$ws.Range("J75").Value = 53
$ws.Range("J76").Value = 258
$ws.Range("J77").Value = 137
$ws.Range("J78").Value = 226
$ws.Range("J79").Value = 512
$ws.Range("J84").Value = 149
$ws.Range("J85").Value = 778
$ws.Range("J90").Value = 205
$ws.Range("J91").Value = 197
$ws.Range("J95").Value = 269
$ws.Range("J99").Value = 266
$ws.Range("I101").Value = 26223
$ws.Range("J101").Value = 17824

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 166
$ws.Range("J7").Value = 557

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J3").Value = 259
$ws.Range("J7").Value = 812

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("H2").Value = 542
$ws.Range("H3").Value = 674
$ws.Range("J3").Value = 350
$ws.Range("J6").Value = 260
$ws.Range("J7").Value = 1012

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J3").Value = 149
$ws.Range("J7").Value = 520

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J2").Value = 43
$ws.Range("J6").Value = 45
$ws.Range("J7").Value = 133

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J6").Value = 101
$ws.Range("J7").Value = 199

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J2").Value = 42
$ws.Range("J3").Value = 53
$ws.Range("J7").Value = 258

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 163
$ws.Range("J7").Value = 515

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("J6").Value = 49
$ws.Range("J7").Value = 133

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("J2").Value = 24
$ws.Range("J7").Value = 115

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 157
$ws.Range("J3").Value = 146
$ws.Range("J6").Value = 368
$ws.Range("J7").Value = 722

$ws = $wb.Worksheets.Item('Boystown')
$ws.Range("J3").Value = 7
$ws.Range("J6").Value = 21

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("J6").Value = 65
$ws.Range("J7").Value = 118

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J2").Value = 60
$ws.Range("J3").Value = 77
$ws.Range("J6").Value = 60
$ws.Range("J7").Value = 226

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J2").Value = 53
$ws.Range("J6").Value = 109

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J2").Value = 45
$ws.Range("J7").Value = 169

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J3").Value = 82
$ws.Range("J7").Value = 197

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J3").Value = 182
$ws.Range("J6").Value = 141
$ws.Range("J7").Value = 512

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("J3").Value = 32
$ws.Range("J7").Value = 119

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J6").Value = 100
$ws.Range("J7").Value = 373

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("J2").Value = 9
$ws.Range("J7").Value = 26

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J4").Value = 11
$ws.Range("J7").Value = 247

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J6").Value = 103
$ws.Range("J7").Value = 272

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("J6").Value = 31
$ws.Range("J7").Value = 105

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("J2").Value = 11
$ws.Range("J7").Value = 57

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J2").Value = 107
$ws.Range("J3").Value = 129
$ws.Range("J7").Value = 451

$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Range("J4").Value = 2
$ws.Range("J7").Value = 28

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("J2").Value = 58
$ws.Range("J3").Value = 46
$ws.Range("J7").Value = 164

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("J2").Value = 38
$ws.Range("J7").Value = 135

$ws = $wb.Worksheets.Item('Galewood')
$ws.Range("J3").Value = 8
$ws.Range("J7").Value = 28

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("J6").Value = 19
$ws.Range("J7").Value = 76

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J3").Value = 20
$ws.Range("J7").Value = 96

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("J6").Value = 11
$ws.Range("J7").Value = 53

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J2").Value = 75
$ws.Range("J7").Value = 205

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J6").Value = 79
$ws.Range("J7").Value = 223

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J2").Value = 65
$ws.Range("J4").Value = 13

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("J3").Value = 49
$ws.Range("J7").Value = 137

$ws = $wb.Worksheets.Item('Andersonville')
$ws.Range("J5").Value = 3
$ws.Range("J6").Value = 7

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J3").Value = 343
$ws.Range("J4").Value = 61
$ws.Range("J6").Value = 381
$ws.Range("J7").Value = 1141
